$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.333.08'
$ws.Range('E2').Value = '  +3.03%  '
$ws.Range('D3').Value = '2.093.11'
$ws.Range('E3').Value = '  +4.49%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '251.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.10%  '
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.16'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +22.41%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '61.76'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.43%  '
$ws.Range('E10').Value = '  +3.40%  '
$ws.Range('E11').Value = '  +4.57%  '
$ws.Range('E12').Value = '  +8.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.14'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.47%  '
$ws.Range('D14').Value = '2.399.42'
$ws.Range('E14').Value = '  +4.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.835'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.27%  '
$ws.Range('D16').Value = '2.100.46'
$ws.Range('E16').Value = '  +4.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.17'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +6.59%  '
$ws.Range('D18').Value = '37.275.54'
$ws.Range('E18').Value = '  +2.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.76'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.56'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +14.86%  '
$ws.Range('E21').Value = '  +5.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '240.78'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.11%  '
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.47'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '172.29'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.30%  '
$ws.Range('E29').Value = '  +3.46%  '
$ws.Range('E30').Value = '  +2.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '23.75'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +12.26%  '
$ws.Range('E32').Value = '  +28.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.61%  '
$ws.Range('E34').Value = '  +6.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0915'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.49%  '
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('E37').Value = '  +5.07%  '
$ws.Range('E38').Value = '  -0.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.26'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.35%  '
$ws.Range('E40').Value = '  +2.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.30'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +15.54%  '
$ws.Range('E42').Value = '  +6.36%  '
$ws.Range('E43').Value = '  +6.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '98.70'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0932'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +15.25%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +111.03%  '
$ws.Range('B47').Value = 'HuobiToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.81'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.08%  '
$ws.Range('D48').Value = '1.322.88'
$ws.Range('E48').Value = '  +1.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.93'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.15%  '
$ws.Range('E50').Value = '  +14.98%  '
$ws.Range('E51').Value = '  +7.05%  '
